$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.434.48"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.580.62"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.68"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.248"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0895"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").Value = "1.808.16"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").Value = "1.576.15"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").Value = "28.438.90"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "230.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "0.0₃0689"
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("D35").Value = "1.399.58"
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("E36").Value = "  +7.19%  "
$ws.Range("E37").Value = "  -4.08%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.523"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.788"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.96%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0457"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.929"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.50%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "62.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "1.719.15"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0518"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "
